# Applies the commit: appends 5 new rows to the single table in the
# document, describing additional Java programs.

$d = $word.ActiveDocument

function Set-CellRuns {
    param($Cell, $Parts, $WithTabs)

    $runsXml = ""
    foreach ($part in $Parts) {
        $escaped = $part.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
        $runsXml += "<w:r><w:t>$escaped</w:t></w:r>"
    }

    $pPrXml = ""
    if ($WithTabs) {
        $pPrXml = "<w:pPr><w:tabs><w:tab w:val=`"left`" w:pos=`"3249`"/></w:tabs></w:pPr>"
    }

    $frag = "<?xml version=`"1.0`" encoding=`"UTF-8`"?><w:wordDocument xmlns:w=`"http://schemas.microsoft.com/office/word/2003/wordml`"><w:body><w:p>$pPrXml$runsXml</w:p></w:body></w:wordDocument>"

    $r = $Cell.Range
    $r.End = $r.End - 1
    $r.InsertXML($frag)
}

# Each entry: Col1 parts (first paragraph keeps the "tabs" pPr seen on the
# other rows of the table), Col2 parts (program file name, no pPr).
$newRows = @(
    @{
        Col1 = @("Reversing a number / finding number is ", "palindrome ", "or not ")
        Col2 = @("PalindromeNumber", ".java")
    },
    @{
        Col1 = @("Reversing a String using StringBuffer")
        Col2 = @("ReverseString1.java")
    },
    @{
        Col1 = @("Reversing String ", "without using reverse function")
        Col2 = @("ReverseString", "2", ".java")
    },
    @{
        Col1 = @("Finding the character occurrence in a string")
        Col2 = @("CharacterOccurance.java")
    },
    @{
        Col1 = @("Sorting the String", " in descending order", " by removing the digits")
        Col2 = @("SortString1.java")
    }
)

foreach ($rowData in $newRows) {
    # Re-fetch the table each time: row collections are considered stale
    # after a structural edit (a row was just added).
    $t = $d.Tables.Item(1)
    $t.Rows.Add() | Out-Null

    $t = $d.Tables.Item(1)
    $newRowIndex = $t.Rows.Count

    $cell1 = $t.Cell($newRowIndex, 1)
    Set-CellRuns $cell1 $rowData.Col1 $true

    $t = $d.Tables.Item(1)
    $cell2 = $t.Cell($newRowIndex, 2)
    Set-CellRuns $cell2 $rowData.Col2 $false
}

$t = $d.Tables.Item(1)
Write-Host "Final row count:" $t.Rows.Count
